$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Hunk 1: "M.Thanes" -> "PreethaRaai", then append a new run " RS"
#         right after the existing <w:proofErr w:type="spellEnd"/>.
# ------------------------------------------------------------------

# 1a. Replace the name text in place (keeps it as a single run so the
#     existing rPr / spellStart / spellEnd wrapping survive untouched).
$r1 = $d.Content
$null = $r1.Find.Execute("M.Thanes", $true, $false, $false, $false, $false, $true, 1, $false, "PreethaRaai", 2)

# 1b. Locate the run we just produced so we know its exact bounds.
$name = $d.Content
$null = $name.Find.Execute("PreethaRaai", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nameStart = $name.Start
$nameEnd = $name.End

# 1c. Insert a placeholder run right after it (InsertAfter always
#     produces an unformatted run, so it is only used to reserve the
#     text/position).
$name.Collapse(0)
$name.InsertAfter(" RS")
$suffixStart = $name.Start
$suffixEnd = $name.End

# 1d. Clone formatting from the (untouched, still valid) "PreethaRaai"
#     run onto the new run via FormattedText, which copies the full
#     rPr faithfully without corrupting neighboring runs.
$nameSrc = $d.Range($nameStart, $nameEnd)
$nameFt = $nameSrc.FormattedText
$suffixDest = $d.Range($suffixStart, $suffixEnd)
$suffixDest.FormattedText = $nameFt

# 1e. FormattedText assignment also copied the source TEXT, so fix the
#     destination text back to " RS" using a position-based range
#     (text-based Find would now be ambiguous).
$fixLen = $nameEnd - $nameStart
$suffixFix = $d.Range($suffixStart, $suffixStart + $fixLen)
$suffixFix.Text = " RS"

# ------------------------------------------------------------------
# Hunk 2: " 2021503712" -> " 2021503", then append a new run "320"
#         at the end of that paragraph.
# ------------------------------------------------------------------

# 2a. Truncate the roll number text in place (single run, rPr intact).
$r2 = $d.Content
$null = $r2.Find.Execute(" 2021503712", $true, $false, $false, $false, $false, $true, 1, $false, " 2021503", 2)

# 2b. Locate the run we just produced.
$roll = $d.Content
$null = $roll.Find.Execute(" 2021503", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rollStart = $roll.Start
$rollEnd = $roll.End

# 2c. Insert a placeholder run right after it.
$roll.Collapse(0)
$roll.InsertAfter("320")
$suffix2Start = $roll.Start
$suffix2End = $roll.End

# 2d. Clone formatting from the roll-number run onto the new run.
$rollSrc = $d.Range($rollStart, $rollEnd)
$rollFt = $rollSrc.FormattedText
$suffix2Dest = $d.Range($suffix2Start, $suffix2End)
$suffix2Dest.FormattedText = $rollFt

# 2e. Fix the destination text back to "320".
$fix2Len = $rollEnd - $rollStart
$suffix2Fix = $d.Range($suffix2Start, $suffix2Start + $fix2Len)
$suffix2Fix.Text = "320"
